$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename the TV-size question to refer to the "primary" TV (since a
#    second TV question block is being introduced below).
$ws.Range("C36").Value = "What is the size of your primary TV?"

# 2. Add a new question block (rows 44-45): "Do you have a secondary TV?"
#    Re-uses the exact formatting of the existing Yes/No question block
#    (rows 14-15, "Are you bilingual?") so styles line up.
$ws.Range("A14:K15").Copy($ws.Range("A44"))
$ws.Range("D44").ClearContents()
$ws.Range("B45:G45").ClearContents()
$ws.Range("K44:K45").ClearContents()

$ws.Range("A44").Value = "11"
$ws.Range("B44").Value = "qx_has_second_tv_key"
$ws.Range("C44").Value = "Do you have a secondary TV?"
$ws.Range("E44").Value = 1
$ws.Range("F44").Value = "2"
$ws.Range("G44").Value = "answer_has_second_tv_key_no"
$ws.Range("H44").Value = "Yes"
$ws.Range("I44").Value = "answer_has_second_tv_key_yes"
$ws.Range("J44").Value = 10

$ws.Range("H45").Value = "No"
$ws.Range("I45").Value = "answer_has_second_tv_key_no"
$ws.Range("J45").Value = 20

# 3. Add a new question block (rows 47-49): "What is the size of your
#    secondary TV?" Re-uses the exact formatting of the existing TV-size
#    question block (rows 36-38).
$ws.Range("A36:K38").Copy($ws.Range("A47"))
$ws.Range("D47").ClearContents()
$ws.Range("B48:G48").ClearContents()
$ws.Range("B49:G49").ClearContents()

$ws.Range("A47").Value = "12"
$ws.Range("B47").Value = "qx_2nd_tv_size_key"
$ws.Range("C47").Value = "What is the size of your secondary TV?"
$ws.Range("F47").Value = "3"
$ws.Range("H47").Value = "32 or less"
$ws.Range("I47").Value = "answer_2nd_tv_size_32_less"
$ws.Range("J47").Value = 10

$ws.Range("H48").Value = "33 to 66"
$ws.Range("I48").Value = "answer_2nd_tv_size_33_66"
$ws.Range("J48").Value = 20

$ws.Range("H49").Value = "66 or less"
$ws.Range("I49").Value = "answer_2nd_tv_size_66_more"
$ws.Range("J49").Value = 30

# 4. Match the author's final view/selection state.
$ws.Range("A33").Select()
$ws.Range("G55").Select()
